# Apply the two edits described by the diff:
#  1. Slide 3 ("Requirements, Goals and Scope"): remove the two bullet
#     paragraphs "High scale for number of test sessions and faster
#     detection interval" (lvl 1) and "Support hardware implementation"
#     (lvl 2) from the Content Placeholder.
#  2. Slide 7 ("Return Path Control Code Sub-TLV - Usage"): change the
#     wording "Session-Reflector sends test packet in-band..." to
#     "Session-Reflector transmits test packet in-band...".

$p = $ppt.ActivePresentation

# --- Slide 3: delete the two paragraphs ---
$s3 = $p.Slides.Item(3)
$shape3 = $s3.Shapes.Item(2)   # "Content Placeholder 2"
$tr3 = $shape3.TextFrame.TextRange

# Delete from the bottom up so earlier paragraph indices stay valid.
$tr3.Paragraphs(7).Delete()
$tr3.Paragraphs(6).Delete()

# --- Slide 7: reword "sends" -> "transmits" ---
$s7 = $p.Slides.Item(7)
$shape7 = $s7.Shapes.Item(2)   # "Content Placeholder 2"
$tr7 = $shape7.TextFrame.TextRange
$para7 = $tr7.Paragraphs(3)
$run7 = $para7.Runs(1)
$run7.Text = "Session-Reflector transmits test packet in-band on the same incoming link in the reverse direction"
